$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3153
$ws1.Range("F3").Value = 731
$ws1.Range("F4").Value = 115
$ws1.Range("F5").Value = 6855
$ws1.Range("F6").Value = 1930
$ws1.Range("F7").Value = 13
$ws1.Range("F11").Value = 65
$ws1.Range("F13").Value = 145
$ws1.Range("F14").Value = 169
$ws1.Range("F15").Value = 32

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 10

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3153
$ws4.Range("F3").Value = 10
$ws4.Range("F4").Value = 731
$ws4.Range("F5").Value = 115
$ws4.Range("F6").Value = 6855
$ws4.Range("F7").Value = 1930
$ws4.Range("F8").Value = 13
$ws4.Range("F12").Value = 65
$ws4.Range("F14").Value = 145
$ws4.Range("F15").Value = 169
$ws4.Range("F16").Value = 32
